# "Updates for my voice" - rephrase several sentences from "we" to "I" style
# narration, add a short aside about printf debugging, add a parenthetical
# remark about platform files, mention the "internal" buffer, and relocate
# the _GoBack bookmark to sit just before the final "application." word.

$d = $word.ActiveDocument

# --- Paragraph: "In this video, we are going to use ..." ---------------
$r = $d.Content
$r.Find.Execute("In this video, we are", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "In this video, I am", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("to a terminal window.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, `
                 "to a terminal window…. You know.. printf debugging that we all like to do", 2) | Out-Null

# --- Paragraph: "We will add information ..." ---------------------------
$r = $d.Content
$r.Find.Execute("We will add information", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "I will add information", 2) | Out-Null

# --- Paragraph: "In 03_blinkled_print.c we are going to add ..." -------
$r = $d.Content
$r.Find.Execute("_blinkled_print.c we are going to add", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "_blinkled_print.c I am going to add", 2) | Out-Null

# --- Paragraph: "The interface is configured and started by default." --
$r = $d.Content
$r.Find.Execute("The interface is configured and started by default.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, `
                 "The interface is configured and started by default (remember those platform file ….)", 2) | Out-Null

# --- Paragraph: "... until the buffer is full." -------------------------
$r = $d.Content
$r.Find.Execute("until the buffer is full.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "until the internal buffer is full.", 2) | Out-Null

# --- Paragraph: "OK, so now that we have added ..." ---------------------
$r = $d.Content
$r.Find.Execute("OK, so now that we have added", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "OK, so now that I have added", 2) | Out-Null

# --- Paragraph: "Once it is programmed, we will open ..." ---------------
$r = $d.Content
$r.Find.Execute("Once it is programmed, we will open", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Once it is programmed, I will open", 2) | Out-Null

# --- Paragraph: "... function that we call at the beginning of our application." --
$r = $d.Content
$r.Find.Execute("function that we call at the beginning of our application.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "function that I called at the beginning of the application.", 2) | Out-Null

# --- Relocate the "_GoBack" bookmark to sit right before "application." --
$srch = $d.Content
$srch.Find.Execute("application.") | Out-Null
$bmRange = $d.Range($srch.Start, $srch.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
